$wb = $excel.ActiveWorkbook

$dateFmt = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------
# Sheet 1: "Weekly Quantity" - append rows 6 and 7
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

$ws1.Cells.Item(6, 1).Value = 45662.99999999999
$ws1.Cells.Item(6, 1).NumberFormat = $dateFmt
$ws1.Cells.Item(6, 2).Value = 1

$ws1.Cells.Item(7, 1).Value = 45669.99999999999
$ws1.Cells.Item(7, 1).NumberFormat = $dateFmt
$ws1.Cells.Item(7, 2).Value = 1

# ---------------------------------------------------------------
# Sheet 2: "Monthly Trend" - append row 6
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Monthly Trend")

$ws2.Cells.Item(6, 1).Value = 45688.99999999999
$ws2.Cells.Item(6, 1).NumberFormat = $dateFmt
$ws2.Cells.Item(6, 2).Value = 2

# ---------------------------------------------------------------
# Sheet 3: "PO Forecast" - new forecast model
# rows 2-3 stay as-is, rows 4-15 get new values
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("PO Forecast")

$poDates = @(
    45613.99999999999,
    45641.99999999999,
    45662.99999999999,
    45669.99999999999,
    45676.99999999999,
    45683.99999999999,
    45690.99999999999,
    45697.99999999999,
    45704.99999999999,
    45711.99999999999,
    45718.99999999999,
    45725.99999999999
)

$row = 4
foreach ($d in $poDates) {
    $ws3.Cells.Item($row, 1).Value = $d
    $ws3.Cells.Item($row, 1).NumberFormat = $dateFmt
    $ws3.Cells.Item($row, 2).Value = 2
    $row = $row + 1
}
